# Fixed issue with ramp constraint. Added initial states for processes and
# for storages (states).
#
# This script applies the following changes to the workbook:
#   1. "nodes" sheet: new column J "initial_state" (all zeros).
#   2. "processes" sheet: new column M "initial_state" (first row = 1, rest 0).
#   3. "process_topology" sheet: ramp limits (columns G/H, rows 2-6) changed
#      to 0.5.
#   4. "inflow" sheet: row 2 (B2:D2 and dependent formulas E2:G2) changed
#      from 0 to -5.
#   5. Selection/active-sheet bookkeeping to match the final workbook state.

$wb = $excel.ActiveWorkbook

# xlCenter
$xlCenter = -4108

# --- 1. "nodes": add "initial_state" column (J) ------------------------
$wsNodes = $wb.Worksheets.Item("nodes")

$wsNodes.Range("J1").Value = "initial_state"

$wsNodes.Range("J2").Value = 0
$wsNodes.Range("J3").Value = 0
$wsNodes.Range("J4").Value = 0
$wsNodes.Range("J5").Value = 0
$wsNodes.Range("J6").Value = 0
$wsNodes.Range("J7").Value = 0
$wsNodes.Range("J2:J7").HorizontalAlignment = $xlCenter

# --- 2. "processes": add "initial_state" column (M) ---------------------
$wsProcesses = $wb.Worksheets.Item("processes")

$wsProcesses.Range("M1").Value = "initial_state"

$wsProcesses.Range("M2").Value = 1
$wsProcesses.Range("M3").Value = 0
$wsProcesses.Range("M4").Value = 0
$wsProcesses.Range("M5").Value = 0
$wsProcesses.Range("M6").Value = 0
$wsProcesses.Range("M7").Value = 0
$wsProcesses.Range("M8").Value = 0
$wsProcesses.Range("M2:M8").HorizontalAlignment = $xlCenter

# --- 3. "process_topology": ramp limits now 0.5 --------------------------
$wsTopology = $wb.Worksheets.Item("process_topology")

$wsTopology.Range("G2").Value = 0.5
$wsTopology.Range("H2").Value = 0.5
$wsTopology.Range("G3").Value = 0.5
$wsTopology.Range("H3").Value = 0.5
$wsTopology.Range("G4").Value = 0.5
$wsTopology.Range("H4").Value = 0.5
$wsTopology.Range("G5").Value = 0.5
$wsTopology.Range("H5").Value = 0.5
$wsTopology.Range("G6").Value = 0.5
$wsTopology.Range("H6").Value = 0.5

# --- 4. "inflow": row 2 values become -5 ---------------------------------
$wsInflow = $wb.Worksheets.Item("inflow")

$wsInflow.Range("B2").Value = -5
$wsInflow.Range("C2").Value = -5
$wsInflow.Range("D2").Value = -5

# --- 5. Selections / active sheet, matching the saved workbook state ----
$wsCf = $wb.Worksheets.Item("cf")

$wsProcesses.Range("M2").Select() | Out-Null
$wsCf.Range("D2").Select() | Out-Null
$wsInflow.Range("B3").Select() | Out-Null
$wsNodes.Range("J4").Select() | Out-Null
$wsTopology.Range("J6").Select() | Out-Null
